$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '26.771.84'
$ws.Cells.Item(2, 5).Value = '  +0.58%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.644.42'
$ws.Cells.Item(3, 5).Value = '  +0.00%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.43%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''216.79'
$ws.Cells.Item(5, 5).Value = '  +0.36%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''0.500'
$ws.Cells.Item(6, 5).Value = '  -0.67%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.33%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.56%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.57%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''19.17'
$ws.Cells.Item(10, 5).Value = '  -0.53%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.15%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.866.38'
$ws.Cells.Item(12, 5).Value = '  -0.41%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.645.14'
$ws.Cells.Item(13, 5).Value = '  +0.04%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''4.17'
$ws.Cells.Item(14, 5).Value = '  -1.31%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -1.30%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''64.33'
$ws.Cells.Item(16, 5).Value = '  -2.73%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '26.767.24'
$ws.Cells.Item(17, 5).Value = '  +0.41%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -1.82%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''214.16'
$ws.Cells.Item(19, 5).Value = '  -1.94%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.35%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''4.37'
$ws.Cells.Item(21, 5).Value = '  -0.31%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +14.32%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.77%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -2.59%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''144.97'
$ws.Cells.Item(25, 5).Value = '  -1.19%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''1.00'
$ws.Cells.Item(26, 5).Value = '  -0.60%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.58%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''7.11'
$ws.Cells.Item(28, 5).Value = '  -0.32%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''15.65'
$ws.Cells.Item(29, 5).Value = '  -1.43%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.52%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.25%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''3.32'
$ws.Cells.Item(32, 5).Value = '  -2.17%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -1.67%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.294.26'
$ws.Cells.Item(34, 5).Value = '  +1.17%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.18%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +1.39%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -5.81%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''0.536'
$ws.Cells.Item(38, 5).Value = '  +1.50%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''0.827'
$ws.Cells.Item(39, 5).Value = '  +0.15%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.34%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.807'
$ws.Cells.Item(41, 5).Value = '  -0.14%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.06%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -1.91%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '1.792.87'
$ws.Cells.Item(44, 5).Value = '  +0.40%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''61.38'
$ws.Cells.Item(45, 5).Value = '  +2.80%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''91.38'
$ws.Cells.Item(46, 5).Value = '  -1.96%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -1.12%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).Value = '''0.0517'
$ws.Cells.Item(48, 5).Value = '  +0.01%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '''7.66'
$ws.Cells.Item(49, 5).Value = '  -2.13%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).Value = '''0.0976'
$ws.Cells.Item(50, 5).Value = '  -0.12%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Mantle'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value = '''0.407'
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
